$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (shifts old N:P -> O:Q), inheriting the
# width of the column to its left (M), as Excel does for a plain column insert.
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and set its selection,
# moving the selection away from the "Transactions" sheet.
$ws.Activate()
$ws.Range("P15").Select()
